# Regenerate the "K" column (column G) save_data values for scott_tanner.xlsx.
# These are recomputed statistics (K values, replacing the old "Strike#"
# derived numbers) that are written back into the sheet as literal values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value, taken from the regenerated
# save_data output.
$kUpdates = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 2
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    19 = 2
    20 = 2
    21 = 2
    22 = 1
    23 = 2
    24 = 1
    25 = 1
    26 = 0
    27 = 0
    28 = 4
    29 = 1
    30 = 2
    31 = 1
    32 = 1
    33 = 2
    34 = 0
    35 = 3
    36 = 3
    37 = 2
    38 = 2
    39 = 1
    40 = 1
    41 = 2
    42 = 1
    43 = 1
    45 = 1
    46 = 2
    47 = 0
    48 = 1
    49 = 2
    50 = 0
    51 = 0
    52 = 1
    53 = 2
    54 = 0
    55 = 1
    56 = 1
    57 = 0
    58 = 0
    59 = 2
    60 = 2
    61 = 0
    62 = 3
    63 = 2
    67 = 2
    68 = 1
    69 = 1
    70 = 2
}

foreach ($row in $kUpdates.Keys) {
    $ws.Cells.Item($row, 7).Value = $kUpdates[$row]
}
